# Generate Report for Handoff
#
# Updates the "Latest Handoff Datetime" for the 23d785a5-... file on the
# zh-cn and de-de status sheets to reflect a new handoff that just ran,
# while leaving every other row/column untouched.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 on each sheet corresponds to file
# 23d785a5-46aa-4652-9eb4-47fb758425bd ; column D is "Latest Handoff Datetime".
$wsZhCn.Range("D4").Value = "2016-01-27 07:33:53"
$wsDeDe.Range("D4").Value = "2016-01-27 07:34:06"
